# Insert two new data rows right before the old row 407 (Ajo / Vega Modelo de
# Temuco, La Araucania), pushing the existing rows 407-425 down to 409-427.
# This mirrors the source diff: dimension A1:R425 -> A1:R427, with two brand
# new rows of weekly price data and everything below shifted down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 407 (each Insert() pushes rows 407+ down by one).
$ws.Rows.Item(407).Insert()
$ws.Rows.Item(407).Insert()

# New row 407
$ws.Range("A407").Value = 10
$ws.Range("B407").Value = "Vega Modelo de Temuco"
$ws.Range("C407").Value = "La Araucanía"
$ws.Range("D407").Value = 44516
$ws.Range("D407").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E407").Value = 9
$ws.Range("F407").Value = 100112003
$ws.Range("G407").Value = "Ajo"
$ws.Range("H407").Value = "Chino"
$ws.Range("I407").Value = "Primera"
$ws.Range("J407").Value = 315
$ws.Range("K407").Value = 19000
$ws.Range("L407").Value = 20000
$ws.Range("M407").Value = 19492
$ws.Range("N407").Value = "`$/caja 10 kilos"
$ws.Range("O407").Value = "China"
$ws.Range("P407").Value = 1949
$ws.Range("Q407").Value = 10
$ws.Range("R407").Value = "Hortaliza"

# New row 408
$ws.Range("A408").Value = 10
$ws.Range("B408").Value = "Vega Modelo de Temuco"
$ws.Range("C408").Value = "La Araucanía"
$ws.Range("D408").Value = 44516
$ws.Range("D408").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E408").Value = 9
$ws.Range("F408").Value = 100112003
$ws.Range("G408").Value = "Ajo"
$ws.Range("H408").Value = "Chino"
$ws.Range("I408").Value = "Primera"
$ws.Range("J408").Value = 65
$ws.Range("K408").Value = 21000
$ws.Range("L408").Value = 21000
$ws.Range("M408").Value = 21000
$ws.Range("N408").Value = "`$/malla 10 kilos"
$ws.Range("O408").Value = "China"
$ws.Range("P408").Value = 2100
$ws.Range("Q408").Value = 10
$ws.Range("R408").Value = "Hortaliza"
